# Apply updated betting-odds values to the "Jogos da Semana" sheet
# (rows 2-8), matching the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.75
$ws.Range("I2").Value = 5
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 7
$ws.Range("Y2").Value = 9.5
$ws.Range("AC2").Value = 6.5
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 23
$ws.Range("AJ2").Value = 51
$ws.Range("AL2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 6.5
$ws.Range("AZ2").Value = 126

# Row 3
$ws.Range("J3").Value = 2.25
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 8.5
$ws.Range("AI3").Value = 19
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 29

# Row 4
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 7
$ws.Range("Z4").Value = 15
$ws.Range("AA4").Value = 21
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AG4").Value = 9.5
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 51
$ws.Range("AK4").Value = 51
$ws.Range("AL4").Value = 67
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 41
$ws.Range("AU4").Value = 10
$ws.Range("AW4").Value = 6.5
$ws.Range("AX4").Value = 34
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 201

# Row 5
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 5.5
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.88
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 12
$ws.Range("AA5").Value = 15
$ws.Range("AD5").Value = 8
$ws.Range("AE5").Value = 19
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 51
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 41
$ws.Range("AN5").Value = 3.6
$ws.Range("AO5").Value = 8.5
$ws.Range("AQ5").Value = 26
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 51
$ws.Range("AW5").Value = 7
$ws.Range("AX5").Value = 29
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 101
$ws.Range("BA5").Value = 126

# Row 7
$ws.Range("G7").Value = 2.2
$ws.Range("I7").Value = 3
$ws.Range("L7").Value = 3.4
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25
$ws.Range("U7").Value = 1.62
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 9.5
$ws.Range("X7").Value = 12
$ws.Range("Z7").Value = 21
$ws.Range("AJ7").Value = 29
$ws.Range("AN7").Value = 4.5
$ws.Range("AO7").Value = 12
$ws.Range("AT7").Value = 3.25
$ws.Range("AW7").Value = 5
$ws.Range("AX7").Value = 15

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("W8").Value = 11
$ws.Range("AC8").Value = 13
$ws.Range("AM8").Value = 151
$ws.Range("AQ8").Value = 51
$ws.Range("BC8").Value = 451
